$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "module" row that was "pwm out mosfet" (0/0) to the new
# "buzzer / pwm / io?" entry with 1 pin / 1 instance.
$ws.Range("E10").Value = "buzzer / pwm / io?"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1

# Re-apply the Total formula across H3:H10 as a shared formula (matches
# how Excel stores a fill-down of the same relative formula).
$ws.Range("H3:H10").Formula = "=F3*G3"

# Move the active selection from E11 to E9.
$ws.Range("E9").Select()
